$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "food_event_fields"

# Update selection
$ws.Range("D8").Select()

# Remove row-level custom format flag (keep per-cell styles identical)
$ws.Rows.Item(1).ClearFormats()
$ws.Range("A1").VerticalAlignment = -4160
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").VerticalAlignment = -4160
$ws.Range("B1").HorizontalAlignment = -4131
$ws.Range("B1").Font.Bold = $true
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Font.Bold = $true

$ws.Rows.Item(2).ClearFormats()
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("D2").VerticalAlignment = -4160

$ws.Rows.Item(3).ClearFormats()
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("B3").VerticalAlignment = -4160
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("D3").VerticalAlignment = -4160

$ws.Rows.Item(4).ClearFormats()
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("B4").VerticalAlignment = -4160
$ws.Range("B4").HorizontalAlignment = -4131
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("D4").VerticalAlignment = -4160

$ws.Rows.Item(5).ClearFormats()
$ws.Range("A5").VerticalAlignment = -4160
$ws.Range("B5").VerticalAlignment = -4160
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("C5").VerticalAlignment = -4160
$ws.Range("D5").VerticalAlignment = -4160
$ws.Range("D5").WrapText = $true

$ws.Rows.Item(6).ClearFormats()
$ws.Range("A6").VerticalAlignment = -4160
$ws.Range("B6").VerticalAlignment = -4160
$ws.Range("B6").HorizontalAlignment = -4131
$ws.Range("C6").VerticalAlignment = -4160
$ws.Range("D6").VerticalAlignment = -4160

$ws.Rows.Item(7).ClearFormats()
$ws.Range("A7").VerticalAlignment = -4160
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("C7").VerticalAlignment = -4160
$ws.Range("D7").VerticalAlignment = -4160
$ws.Range("D7").WrapText = $true

$ws.Rows.Item(8).ClearFormats()
$ws.Range("A8").VerticalAlignment = -4160
$ws.Range("B8").VerticalAlignment = -4160
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("C8").VerticalAlignment = -4160
$ws.Range("D8").VerticalAlignment = -4160
$ws.Range("D8").WrapText = $true

$ws.Rows.Item(9).ClearFormats()
$ws.Range("A9").VerticalAlignment = -4160
$ws.Range("B9").VerticalAlignment = -4160
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("C9").VerticalAlignment = -4160
$ws.Range("D9").VerticalAlignment = -4160
$ws.Range("D9").WrapText = $true

$ws.Rows.Item(10).ClearFormats()
$ws.Range("A10").VerticalAlignment = -4160
$ws.Range("B10").VerticalAlignment = -4160
$ws.Range("B10").HorizontalAlignment = -4131
$ws.Range("C10").VerticalAlignment = -4160
$ws.Range("D10").VerticalAlignment = -4160
$ws.Range("D10").WrapText = $true

$ws.Rows.Item(11).ClearFormats()
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("C11").VerticalAlignment = -4160
$ws.Range("D11").VerticalAlignment = -4160
$ws.Range("D11").WrapText = $true

$ws.Rows.Item(12).ClearFormats()
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("C12").VerticalAlignment = -4160
$ws.Range("D12").VerticalAlignment = -4160
$ws.Range("D12").WrapText = $true

